# Add a new "Magnesiumchlorid 0,5 molar B. Braun" package row (row 17) to
# Sheet1, right after the existing last data row (row 16), copying row 16's
# formatting and filling in the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 16 (formats/styles) into a fresh row 17.
$ws.Rows("16:16").Copy()
$ws.Rows("17:17").Insert()
$ws.Rows("17:17").RowHeight = 12.75

# Populate the new row with the new package's data.
$ws.Range("A17").Value = 45882
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Magnesiumchlorid 0,5 molar B. Braun, Zusatzampulle für              Infusionslösungen"
$ws.Range("D17").Value = "B. Braun Medical AG"
$ws.Range("E17").Value = "05.03.2."
$ws.Range("F17").Value = "B05XA11"
$ws.Range("G17").Value = "Synthetika human"
$ws.Range("H17").Value = 31180
$ws.Range("I17").Value = 31180
$ws.Range("J17").Value = 43340
$ws.Range("K17").Value = 20
$ws.Range("L17").Value = "5 x 10 mL"
$ws.Range("M17").Value = "Ampulle(n)"
$ws.Range("N17").Value = "B"
$ws.Range("O17").Value = "magnesium, chloridum"
$ws.Range("P17").Value = "magnesium 500 mmol, chloridum 1000 mmol, aqua ad iniectabilia q.s. ad solutionem pro 1000 ml."
$ws.Range("Q17").Value = "Magnesiummangel"

# Match the author's final selection (clicked into the new row's name cell).
$ws.Range("C17").Select()
